# Atualização automática via cronjob
# Remove the oldest day's records (2025-04-28, originally rows 2-6) and
# refresh the "estoque_atualizado" (G) and id (A) values for the
# remaining rows, mirroring a scheduled data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the five rows belonging to 2025-04-28 and shift the rows below
# them upward so the sheet only keeps 2025-04-30 through 2025-05-09.
$ws.Range("A2:H6").Delete(-4162)

# After the shift, update the "id" (A) and "estoque_atualizado" (G)
# columns with their refreshed values for each remaining record.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 7).Value = 863

$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 7).Value = 183

$ws.Cells.Item(4, 1).Value = 0
$ws.Cells.Item(4, 7).Value = 3

$ws.Cells.Item(5, 1).Value = 6
$ws.Cells.Item(5, 7).Value = -14

$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(6, 7).Value = 25

$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 7).Value = -3

$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(8, 7).Value = -2

$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 7).Value = -2

$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 7).Value = 23
